# This edit swaps the full data of row 5 and row 6 in the sheet
# (all columns that differ between the two rows are exchanged; columns
# that already held identical values in both rows are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellValue {
    param($cell, $value)
    $cell.Value = $value
}

function Set-CellEmptyPresent {
    param($cell)
    $cell.Value = ""
    $b = $cell.Font.Bold
    $cell.Font.Bold = $b
}

function Set-CellAbsent {
    param($cell)
    $cell.ClearContents()
}

# --- Row 5 (becomes old row6 values) ---
Set-CellValue $ws.Cells.Item(5,1) 131067826
Set-CellValue $ws.Cells.Item(5,2) 79243
Set-CellValue $ws.Cells.Item(5,5) 6425
Set-CellValue $ws.Cells.Item(5,6) "Garnlav"
Set-CellValue $ws.Cells.Item(5,7) "Alectoria sarmentosa"
Set-CellValue $ws.Cells.Item(5,8) "(Ach.) Ach."
Set-CellEmptyPresent $ws.Cells.Item(5,10)
Set-CellAbsent $ws.Cells.Item(5,12)
Set-CellAbsent $ws.Cells.Item(5,13)
Set-CellValue $ws.Cells.Item(5,16) "Långan Öst, Jmt"
Set-CellValue $ws.Cells.Item(5,17) 465891
Set-CellValue $ws.Cells.Item(5,18) 7046290
Set-CellAbsent $ws.Cells.Item(5,26)
Set-CellAbsent $ws.Cells.Item(5,28)
Set-CellAbsent $ws.Cells.Item(5,29)
Set-CellEmptyPresent $ws.Cells.Item(5,32)
Set-CellValue $ws.Cells.Item(5,34) "Granskog"
Set-CellValue $ws.Cells.Item(5,49) "Kristian Zackrisson"
Set-CellValue $ws.Cells.Item(5,50) "Kristian Zackrisson"

# --- Row 6 (becomes old row5 values) ---
Set-CellValue $ws.Cells.Item(6,1) 131067473
Set-CellValue $ws.Cells.Item(6,2) 57884
Set-CellValue $ws.Cells.Item(6,5) 100109
Set-CellValue $ws.Cells.Item(6,6) "Tretåig hackspett"
Set-CellValue $ws.Cells.Item(6,7) "Picoides tridactylus"
Set-CellValue $ws.Cells.Item(6,8) "(Linnaeus, 1758)"
Set-CellAbsent $ws.Cells.Item(6,10)
Set-CellEmptyPresent $ws.Cells.Item(6,12)
Set-CellValue $ws.Cells.Item(6,13) "färska spår"
Set-CellValue $ws.Cells.Item(6,16) "Åbogen, Jmt"
Set-CellValue $ws.Cells.Item(6,17) 465809
Set-CellValue $ws.Cells.Item(6,18) 7046259
Set-CellValue $ws.Cells.Item(6,26) "15:46"
Set-CellValue $ws.Cells.Item(6,28) "15:46"
Set-CellValue $ws.Cells.Item(6,29) "Färska ringhack"
Set-CellAbsent $ws.Cells.Item(6,32)
Set-CellAbsent $ws.Cells.Item(6,34)
Set-CellValue $ws.Cells.Item(6,49) "Elin Albrechtsson"
Set-CellValue $ws.Cells.Item(6,50) "Elin Albrechtsson"
